# Apply fitness value corrections to column C (rows 2-139) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11939
$ws.Range("C3").Value = 10907
$ws.Range("C4").Value = 10907
$ws.Range("C5").Value = 10665
$ws.Range("C6").Value = 10665
$ws.Range("C7").Value = 10401
$ws.Range("C8").Value = 10401
$ws.Range("C9").Value = 10401
$ws.Range("C10").Value = 10401
$ws.Range("C11").Value = 10339
$ws.Range("C12").Value = 10339
$ws.Range("C13").Value = 9805
$ws.Range("C14").Value = 9805
$ws.Range("C15").Value = 9805
$ws.Range("C16").Value = 9805
$ws.Range("C17").Value = 9805
$ws.Range("C18").Value = 9805
$ws.Range("C19").Value = 9805
$ws.Range("C20").Value = 9805
$ws.Range("C21").Value = 9283
$ws.Range("C22").Value = 9283
$ws.Range("C23").Value = 9113
$ws.Range("C24").Value = 9113
$ws.Range("C25").Value = 8838
$ws.Range("C26").Value = 8838
$ws.Range("C27").Value = 8838
$ws.Range("C28").Value = 8771
$ws.Range("C29").Value = 8771
$ws.Range("C30").Value = 8771
$ws.Range("C31").Value = 8771
$ws.Range("C32").Value = 8771
$ws.Range("C33").Value = 8451
$ws.Range("C34").Value = 8027
$ws.Range("C35").Value = 8027
$ws.Range("C36").Value = 8027
$ws.Range("C37").Value = 8027
$ws.Range("C38").Value = 8027
$ws.Range("C39").Value = 8027
$ws.Range("C40").Value = 8027
$ws.Range("C41").Value = 8027
$ws.Range("C42").Value = 8027
$ws.Range("C43").Value = 8027
$ws.Range("C44").Value = 8027
$ws.Range("C45").Value = 8027
$ws.Range("C46").Value = 8004
$ws.Range("C47").Value = 8004
$ws.Range("C48").Value = 8004
$ws.Range("C49").Value = 8004
$ws.Range("C50").Value = 8004
$ws.Range("C51").Value = 8004
$ws.Range("C52").Value = 8004
$ws.Range("C53").Value = 8004
$ws.Range("C54").Value = 7861
$ws.Range("C55").Value = 7861
$ws.Range("C56").Value = 7619
$ws.Range("C57").Value = 7619
$ws.Range("C58").Value = 7619
$ws.Range("C59").Value = 7619
$ws.Range("C60").Value = 7619
$ws.Range("C61").Value = 7619
$ws.Range("C62").Value = 7619
$ws.Range("C63").Value = 7619
$ws.Range("C64").Value = 7619
$ws.Range("C65").Value = 7619
$ws.Range("C66").Value = 7569
$ws.Range("C67").Value = 7569
$ws.Range("C68").Value = 7569
$ws.Range("C69").Value = 7569
$ws.Range("C70").Value = 7569
$ws.Range("C71").Value = 7569
$ws.Range("C72").Value = 7569
$ws.Range("C73").Value = 7569
$ws.Range("C74").Value = 7569
$ws.Range("C75").Value = 7569
$ws.Range("C76").Value = 7569
$ws.Range("C77").Value = 7569
$ws.Range("C78").Value = 7569
$ws.Range("C79").Value = 7569
$ws.Range("C80").Value = 7569
$ws.Range("C81").Value = 7569
$ws.Range("C82").Value = 7569
$ws.Range("C83").Value = 7569
$ws.Range("C84").Value = 7569
$ws.Range("C85").Value = 7569
$ws.Range("C86").Value = 7569
$ws.Range("C87").Value = 7569
$ws.Range("C88").Value = 7569
$ws.Range("C89").Value = 7569
$ws.Range("C90").Value = 7569
$ws.Range("C91").Value = 7569
$ws.Range("C92").Value = 7569
$ws.Range("C93").Value = 7569
$ws.Range("C94").Value = 7569
$ws.Range("C95").Value = 7569
$ws.Range("C96").Value = 7569
$ws.Range("C97").Value = 7569
$ws.Range("C98").Value = 7569
$ws.Range("C99").Value = 7569
$ws.Range("C100").Value = 7569
$ws.Range("C101").Value = 7569
$ws.Range("C102").Value = 7569
$ws.Range("C103").Value = 7569
$ws.Range("C104").Value = 7569
$ws.Range("C105").Value = 7569
$ws.Range("C106").Value = 7569
$ws.Range("C107").Value = 7569
$ws.Range("C108").Value = 7569
$ws.Range("C109").Value = 7569
$ws.Range("C110").Value = 7569
$ws.Range("C111").Value = 7569
$ws.Range("C112").Value = 7569
$ws.Range("C113").Value = 7569
$ws.Range("C114").Value = 7569
$ws.Range("C115").Value = 7569
$ws.Range("C116").Value = 7569
$ws.Range("C117").Value = 7569
$ws.Range("C118").Value = 7569
$ws.Range("C119").Value = 7569
$ws.Range("C120").Value = 7569
$ws.Range("C121").Value = 7569
$ws.Range("C122").Value = 7569
$ws.Range("C123").Value = 7569
$ws.Range("C124").Value = 7569
$ws.Range("C125").Value = 7569
$ws.Range("C126").Value = 7569
$ws.Range("C127").Value = 7569
$ws.Range("C128").Value = 7569
$ws.Range("C129").Value = 7569
$ws.Range("C130").Value = 7569
$ws.Range("C131").Value = 7569
$ws.Range("C132").Value = 7569
$ws.Range("C133").Value = 7569
$ws.Range("C134").Value = 7569
$ws.Range("C135").Value = 7569
$ws.Range("C136").Value = 7569
$ws.Range("C137").Value = 7569
$ws.Range("C138").Value = 7569
$ws.Range("C139").Value = 7569
